# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, text first ...
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# ... then copy the existing header formatting (bold font, borders,
# center/top alignment) from an existing header cell so the new headers
# reuse the same style record instead of minting a new one.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Season record is the same for every player on the roster (team-level
# stat repeated down the column): 68 wins, 94 losses, 0 ties.
$wins = 68
$losses = 94
$ties = 0

$lastRow = 45
$ws.Range("AD2:AD$lastRow").Value = $wins
$ws.Range("AE2:AE$lastRow").Value = $losses
$ws.Range("AF2:AF$lastRow").Value = $ties
